$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 2400
$ws.Range("I7").Value = 2100
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2100
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -1988
$ws.Range("N7").Value = -3224
$ws.Range("H9").Value = 725.75
$ws.Range("I9").Value = 800
$ws.Range("J9").Value = 701
$ws.Range("K9").Value = 800
$ws.Range("L9").Value = 701
$ws.Range("M9").Value = -631
$ws.Range("N9").Value = -1039
$ws.Range("H14").Value = 2400
$ws.Range("I14").Value = 2100
$ws.Range("J14").Value = 3000
$ws.Range("K14").Value = 2100
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = -1909
$ws.Range("N14").Value = -3382
$ws.Range("H18").Value = 16384.125
$ws.Range("I18").Value = 20331.834
$ws.Range("K18").Value = 20331.834
$ws.Range("M18").Value = -20047.834
$ws.Range("H32").Value = 3442.6
$ws.Range("I32").Value = 4000
$ws.Range("K32").Value = 4000
$ws.Range("M32").Value = -3674
$ws.Range("H62").Value = 7247.25
$ws.Range("J62").Value = 11995
$ws.Range("L62").Value = 11995
$ws.Range("N62").Value = -13243
$ws.Range("H65").Value = 7247.25
$ws.Range("J65").Value = 11995
$ws.Range("L65").Value = 59975
$ws.Range("N65").Value = -66215
$ws.Range("H80").Value = 1582.7407
$ws.Range("I80").Value = 910.8
$ws.Range("J80").Value = 1978
$ws.Range("K80").Value = 2732.4
$ws.Range("L80").Value = 5934
$ws.Range("M80").Value = -1734.4
$ws.Range("N80").Value = -7930
$ws.Range("H83").Value = 1582.7407
$ws.Range("I83").Value = 910.8
$ws.Range("J83").Value = 1978
$ws.Range("K83").Value = 8197.199999999999
$ws.Range("L83").Value = 17802
$ws.Range("M83").Value = -3205.199999999999
$ws.Range("N83").Value = -27786
$ws.Range("H106").Value = 11890.637
$ws.Range("I106").Value = 2755.7778
$ws.Range("K106").Value = 2755.7778
$ws.Range("M106").Value = -2124.7778
$ws.Range("H125").Value = 1539.6471
$ws.Range("I125").Value = 877.5
$ws.Range("K125").Value = 7897.5
$ws.Range("M125").Value = -5437.5
$ws.Range("H137").Value = 1155.1
$ws.Range("I137").Value = 1074
$ws.Range("K137").Value = 3222
$ws.Range("M137").Value = -672
$ws.Range("H138").Value = 2111.077
$ws.Range("J138").Value = 2240.625
$ws.Range("L138").Value = 6721.875
$ws.Range("N138").Value = -17001.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6426.227
$ws.Range("I32").Value = 2597.3906
$ws.Range("K32").Value = 2597.3906
$ws.Range("M32").Value = -2310.3906
$ws.Range("H55").Value = 19999.857
$ws.Range("J55").Value = 19999.857
$ws.Range("L55").Value = 19999.857
$ws.Range("N55").Value = -20629.857
$ws.Range("H61").Value = 2390
$ws.Range("I61").Value = 2390
$ws.Range("K61").Value = 2390
$ws.Range("M61").Value = -2178
$ws.Range("H74").Value = 1751.2307
$ws.Range("I74").Value = 1708.3636
$ws.Range("J74").Value = 1987
$ws.Range("K74").Value = 1708.3636
$ws.Range("L74").Value = 1987
$ws.Range("M74").Value = -834.3635999999999
$ws.Range("N74").Value = -3735
$ws.Range("H77").Value = 1751.2307
$ws.Range("I77").Value = 1708.3636
$ws.Range("J77").Value = 1987
$ws.Range("K77").Value = 8541.817999999999
$ws.Range("L77").Value = 9935
$ws.Range("M77").Value = -4173.817999999999
$ws.Range("N77").Value = -18671
$ws.Range("H102").Value = 1655.931
$ws.Range("I102").Value = 1611.7693
$ws.Range("K102").Value = 1611.7693
$ws.Range("M102").Value = 10.23070000000007
$ws.Range("H122").Value = 1947.1875
$ws.Range("I122").Value = 1801.7693
$ws.Range("K122").Value = 5405.3079
$ws.Range("M122").Value = -2955.3079
$ws.Range("H136").Value = 2390
$ws.Range("I136").Value = 2390
$ws.Range("K136").Value = 7170
$ws.Range("M136").Value = -4620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 13102.7
$ws.Range("I105").Value = 16861.143
$ws.Range("K105").Value = 16861.143
$ws.Range("M105").Value = -15114.143
$ws.Range("H141").Value = 58195
$ws.Range("J141").Value = 58195
$ws.Range("L141").Value = 58195
$ws.Range("N141").Value = -68555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 877.3077
$ws.Range("I58").Value = 881
$ws.Range("J58").Value = 857
$ws.Range("K58").Value = 881
$ws.Range("L58").Value = 857
$ws.Range("M58").Value = -678
$ws.Range("N58").Value = -1263
$ws.Range("H62").Value = 3699.6
$ws.Range("I62").Value = 3924.5
$ws.Range("K62").Value = 3924.5
$ws.Range("M62").Value = -3300.5
$ws.Range("H65").Value = 3699.6
$ws.Range("I65").Value = 3924.5
$ws.Range("K65").Value = 19622.5
$ws.Range("M65").Value = -16502.5
$ws.Range("H107").Value = 2089.6191
$ws.Range("I107").Value = 1880.6666
$ws.Range("J107").Value = 2612
$ws.Range("K107").Value = 1880.6666
$ws.Range("L107").Value = 2612
$ws.Range("M107").Value = 39.33339999999998
$ws.Range("N107").Value = -6452
$ws.Range("H135").Value = 106389
$ws.Range("J135").Value = 106389
$ws.Range("L135").Value = 106389
$ws.Range("N135").Value = -116529
$ws.Range("H136").Value = 877.3077
$ws.Range("I136").Value = 881
$ws.Range("J136").Value = 857
$ws.Range("K136").Value = 2643
$ws.Range("L136").Value = 2571
$ws.Range("M136").Value = -93
$ws.Range("N136").Value = -7671
$ws.Range("H141").Value = 239619.81
$ws.Range("J141").Value = 239619.81
$ws.Range("L141").Value = 239619.81
$ws.Range("N141").Value = -249979.81

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 704
$ws.Range("I5").Value = 704
$ws.Range("K5").Value = 2112
$ws.Range("M5").Value = -2000
$ws.Range("H12").Value = 231.25
$ws.Range("I12").Value = 140.33333
$ws.Range("J12").Value = 285.8
$ws.Range("K12").Value = 420.99999
$ws.Range("L12").Value = 857.4000000000001
$ws.Range("M12").Value = -247.99999
$ws.Range("N12").Value = -1203.4
$ws.Range("H88").Value = 9919.429
$ws.Range("J88").Value = 9919.429
$ws.Range("L88").Value = 29758.287
$ws.Range("N88").Value = -30614.287
$ws.Range("H91").Value = 9919.429
$ws.Range("J91").Value = 9919.429
$ws.Range("L91").Value = 29758.287
$ws.Range("N91").Value = -32722.287
$ws.Range("H106").Value = 5801.857
$ws.Range("J106").Value = 5892
$ws.Range("L106").Value = 17676
$ws.Range("N106").Value = -19568
$ws.Range("H120").Value = 6385
$ws.Range("I120").Value = 5023.5
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 15070.5
$ws.Range("L120").Value = 60000
$ws.Range("M120").Value = -10232.5
$ws.Range("N120").Value = -69676
$ws.Range("H122").Value = 885.6
$ws.Range("I122").Value = 799.5
$ws.Range("J122").Value = 943
$ws.Range("K122").Value = 7195.5
$ws.Range("L122").Value = 8487
$ws.Range("M122").Value = -4745.5
$ws.Range("N122").Value = -13387
$ws.Range("H129").Value = 3759.5
$ws.Range("I129").Value = 1225
$ws.Range("J129").Value = 5449.1665
$ws.Range("K129").Value = 3675
$ws.Range("L129").Value = 16347.4995
$ws.Range("M129").Value = 1325
$ws.Range("N129").Value = -26347.4995
$ws.Range("H135").Value = 704
$ws.Range("I135").Value = 704
$ws.Range("K135").Value = 6336
$ws.Range("M135").Value = -3801

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4101.467
$ws.Range("I80").Value = 3199.9
$ws.Range("J80").Value = 5904.6
$ws.Range("K80").Value = 3199.9
$ws.Range("L80").Value = 5904.6
$ws.Range("M80").Value = -2201.9
$ws.Range("N80").Value = -7900.6
$ws.Range("H83").Value = 4101.467
$ws.Range("I83").Value = 3199.9
$ws.Range("J83").Value = 5904.6
$ws.Range("K83").Value = 15999.5
$ws.Range("L83").Value = 29523
$ws.Range("M83").Value = -11007.5
$ws.Range("N83").Value = -39507
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744
$ws.Range("H97").Value = 42890
$ws.Range("I97").Value = 26676.08
$ws.Range("J97").Value = 110448
$ws.Range("K97").Value = 26676.08
$ws.Range("L97").Value = 110448
$ws.Range("M97").Value = -26180.08
$ws.Range("N97").Value = -111440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 43512.11
$ws.Range("I100").Value = 8995
$ws.Range("K100").Value = 8995
$ws.Range("M100").Value = -8454
$ws.Range("H121").Value = 110420
$ws.Range("J121").Value = 110420
$ws.Range("L121").Value = 110420
$ws.Range("N121").Value = -113914
$ws.Range("H132").Value = 3504.639
$ws.Range("I132").Value = 3390.5518
$ws.Range("J132").Value = 3977.2856
$ws.Range("K132").Value = 10171.6554
$ws.Range("L132").Value = 11931.8568
$ws.Range("M132").Value = -7641.6554
$ws.Range("N132").Value = -16991.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 100
$ws.Range("K11").Value = 100
$ws.Range("M11").Value = 42
$ws.Range("H136").Value = 3662.0557
$ws.Range("I136").Value = 2932.3125
$ws.Range("K136").Value = 8796.9375
$ws.Range("M136").Value = -6246.9375
$ws.Range("H138").Value = 95000
$ws.Range("J138").Value = 95000
$ws.Range("L138").Value = 95000
$ws.Range("N138").Value = -105280

Write-Output "Applied 250 cell updates"